$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetDims = @(1, 2, 3, 5, 6, 7, 8, 10)
$stress = @(
    0.5096280925201884,
    0.367408146759599,
    0.3138658379593925,
    0.2905986517730177,
    0.2879770826018936,
    0.3142901195739835,
    0.3184148019222885,
    0.3431524782173103
)

for ($i = 0; $i -lt $targetDims.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $targetDims[$i]
    $ws.Cells.Item($row, 2).Value = $stress[$i]
}
